# Corrected wrong buttons polarity.
# Fix up the BOM rows for the SMD resistors/capacitor (R1,R2 / R3,R4 / C1,C2)
# that previously held the wrong part values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: was mislabeled as resistor "R2"; fix the part designator to C1, C2 capacitor
$ws.Range("D5").Value = "C1, C2"

# Row 3: R1, R2 resistor -- value + MPN were wrong
$ws.Range("B3").Value = "2.2k 100mW"
$ws.Range("E3").Value = "ERJ3GEYJ222V"

# Row 4: R3, R4 resistor -- value + MPN were wrong
$ws.Range("B4").Value = "4.3k 100mW"
$ws.Range("E4").Value = "ERJ3GEYJ432V"

# Row 5: value + MPN were wrong
$ws.Range("B5").Value = "1uF 25V"
$ws.Range("E5").Value = "CL10A105KA8NNNC"

# Restore the active cell selection to B4
$ws.Range("B4").Select()
